# Update the worksheet date and every "axb=c" answer cell to the new
# values from the commit. Each call uses MatchCase + MatchWholeWord so a
# replacement text can never be re-matched by a later Find (the one
# pair that would otherwise collide -- "56x57=3192"->"83x71=5893" and
# the original "83x71=5893"->"93x49=4557" -- is ordered so the latter
# runs first, consuming the original occurrence before the former's
# replacement text exists in the document).
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-01 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-02 Saturday", 2) | Out-Null
$d.Content.Find.Execute("55×64=3520", $true, $false, $false, $false, $false, $true, 1, $false, "71×90=6390", 2) | Out-Null
$d.Content.Find.Execute("39×60=2340", $true, $false, $false, $false, $false, $true, 1, $false, "39×69=2691", 2) | Out-Null
$d.Content.Find.Execute("13×29=377", $true, $false, $false, $false, $false, $true, 1, $false, "97×62=6014", 2) | Out-Null
$d.Content.Find.Execute("85×75=6375", $true, $false, $false, $false, $false, $true, 1, $false, "81×76=6156", 2) | Out-Null
$d.Content.Find.Execute("89×41=3649", $true, $false, $false, $false, $false, $true, 1, $false, "16×71=1136", 2) | Out-Null
$d.Content.Find.Execute("87×74=6438", $true, $false, $false, $false, $false, $true, 1, $false, "89×71=6319", 2) | Out-Null
$d.Content.Find.Execute("51×33=1683", $true, $false, $false, $false, $false, $true, 1, $false, "45×24=1080", 2) | Out-Null
$d.Content.Find.Execute("69×29=2001", $true, $false, $false, $false, $false, $true, 1, $false, "42×61=2562", 2) | Out-Null
$d.Content.Find.Execute("33×96=3168", $true, $false, $false, $false, $false, $true, 1, $false, "76×19=1444", 2) | Out-Null
$d.Content.Find.Execute("56×25=1400", $true, $false, $false, $false, $false, $true, 1, $false, "58×23=1334", 2) | Out-Null
$d.Content.Find.Execute("70×28=1960", $true, $false, $false, $false, $false, $true, 1, $false, "42×51=2142", 2) | Out-Null
$d.Content.Find.Execute("70×49=3430", $true, $false, $false, $false, $false, $true, 1, $false, "63×96=6048", 2) | Out-Null
$d.Content.Find.Execute("71×94=6674", $true, $false, $false, $false, $false, $true, 1, $false, "14×72=1008", 2) | Out-Null
$d.Content.Find.Execute("41×39=1599", $true, $false, $false, $false, $false, $true, 1, $false, "86×37=3182", 2) | Out-Null
$d.Content.Find.Execute("74×28=2072", $true, $false, $false, $false, $false, $true, 1, $false, "88×46=4048", 2) | Out-Null
$d.Content.Find.Execute("65×99=6435", $true, $false, $false, $false, $false, $true, 1, $false, "66×61=4026", 2) | Out-Null
$d.Content.Find.Execute("31×28=868", $true, $false, $false, $false, $false, $true, 1, $false, "40×16=640", 2) | Out-Null
$d.Content.Find.Execute("55×72=3960", $true, $false, $false, $false, $false, $true, 1, $false, "18×54=972", 2) | Out-Null
$d.Content.Find.Execute("62×15=930", $true, $false, $false, $false, $false, $true, 1, $false, "34×32=1088", 2) | Out-Null
$d.Content.Find.Execute("83×71=5893", $true, $false, $false, $false, $false, $true, 1, $false, "93×49=4557", 2) | Out-Null
$d.Content.Find.Execute("56×57=3192", $true, $false, $false, $false, $false, $true, 1, $false, "83×71=5893", 2) | Out-Null
$d.Content.Find.Execute("21×32=672", $true, $false, $false, $false, $false, $true, 1, $false, "52×85=4420", 2) | Out-Null
$d.Content.Find.Execute("24×44=1056", $true, $false, $false, $false, $false, $true, 1, $false, "73×31=2263", 2) | Out-Null
$d.Content.Find.Execute("11×62=682", $true, $false, $false, $false, $false, $true, 1, $false, "27×65=1755", 2) | Out-Null
$d.Content.Find.Execute("64×70=4480", $true, $false, $false, $false, $false, $true, 1, $false, "75×72=5400", 2) | Out-Null
